# Updates odds values on Sheet1 to reflect the latest FlashScore data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("Q2").Value = 1.44
$ws.Range("R2").Value = 2.7

# Row 3
$ws.Range("I3").Value = 5.75
$ws.Range("J3").Value = 2.4
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("U3").Value = 2.38
$ws.Range("V3").Value = 1.53
$ws.Range("X3").Value = 6.5
$ws.Range("Z3").Value = 12
$ws.Range("AK3").Value = 67
$ws.Range("AN3").Value = 3.4
$ws.Range("AW3").Value = 7
$ws.Range("BA3").Value = 201

# Row 4
$ws.Range("M4").Value = 1.17
$ws.Range("N4").Value = 5
$ws.Range("T4").Value = 2

# Row 5
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 8

# Row 6
$ws.Range("G6").Value = 2.15
$ws.Range("H6").Value = 3.4
$ws.Range("I6").Value = 3.4
$ws.Range("J6").Value = 2.63
$ws.Range("K6").Value = 2.3
$ws.Range("L6").Value = 3.6
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 13
$ws.Range("Q6").Value = 1.7
$ws.Range("R6").Value = 2.1
$ws.Range("S6").Value = 1.33
$ws.Range("T6").Value = 3.25
$ws.Range("X6").Value = 12
$ws.Range("Z6").Value = 21
$ws.Range("AC6").Value = 13
$ws.Range("AD6").Value = 6.5
$ws.Range("AE6").Value = 11
$ws.Range("AI6").Value = 19
$ws.Range("AK6").Value = 34
$ws.Range("AL6").Value = 23
$ws.Range("AP6").Value = 19
$ws.Range("AT6").Value = 3.25
